$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 30
$ws.Range("B2").Value = " ack flag count"

$ws.Range("A3").Value = 35
$ws.Range("B3").Value = " active max"

$ws.Range("A4").Value = 36
$ws.Range("B4").Value = " active min"

$ws.Range("A5").Value = 34
$ws.Range("B5").Value = " active std"

$ws.Range("A6").Value = 23
$ws.Range("B6").Value = " bwd iat max"

$ws.Range("A7").Value = 21
$ws.Range("B7").Value = " bwd iat mean"

$ws.Range("A8").Value = 22
$ws.Range("B8").Value = " bwd iat std"

$ws.Range("A9").Value = 26
$ws.Range("B9").Value = " bwd packets/s"

$ws.Range("A10").Value = 2
$ws.Range("B10").Value = " destination port"

$ws.Range("A18").Value = 31

$ws.Range("A23").Value = 9
$ws.Range("B23").Value = " fwd packet length max"

$ws.Range("A24").Value = 10

$ws.Range("A27").Value = 27

$ws.Range("A28").Value = 3
$ws.Range("B28").Value = " protocol"

$ws.Range("A29").Value = 1
$ws.Range("B29").Value = " source port"

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = " syn flag count"

$ws.Range("A34").Value = 33

$ws.Range("A35").Value = 28
$ws.Range("B35").Value = "fin flag count"

$ws.Range("A37").Value = 25

$ws.Range("A39").Value = 32
